$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "89.234.33"
$ws.Cells.Item(2, 5).Value = "  +3.30%  "
$ws.Cells.Item(3, 4).Value = "3.221.62"
$ws.Cells.Item(3, 5).Value = "  -1.53%  "
$ws.Cells.Item(4, 4).Value = "'0.998"
$ws.Cells.Item(4, 5).Value = "  -0.32%  "
$ws.Cells.Item(5, 4).Value = "'212.39"
$ws.Cells.Item(5, 5).Value = "  +0.96%  "
$ws.Cells.Item(6, 4).Value = "'623.06"
$ws.Cells.Item(6, 5).Value = "  -0.07%  "
$ws.Cells.Item(7, 4).Value = "'0.397"
$ws.Cells.Item(7, 5).Value = "  +8.31%  "
$ws.Cells.Item(8, 4).Value = "'0.700"
$ws.Cells.Item(8, 5).Value = "  +8.30%  "
$ws.Cells.Item(9, 4).Value = "'0.997"
$ws.Cells.Item(10, 4).Value = "3.203.29"
$ws.Cells.Item(10, 5).Value = "  -2.30%  "
$ws.Cells.Item(11, 4).Value = "'0.563"
$ws.Cells.Item(11, 5).Value = "  -1.82%  "
$ws.Cells.Item(12, 5).Value = "  +1.42%  "
$ws.Cells.Item(13, 4).Value = "'0.0000259"
$ws.Cells.Item(13, 5).Value = "  +1.02%  "
$ws.Cells.Item(14, 4).Value = "'5.39"
$ws.Cells.Item(14, 5).Value = "  +2.38%  "
$ws.Cells.Item(15, 4).Value = "3.807.39"
$ws.Cells.Item(15, 5).Value = "  -1.62%  "
$ws.Cells.Item(16, 4).Value = "'33.45"
$ws.Cells.Item(16, 5).Value = "  -1.42%  "
$ws.Cells.Item(17, 4).Value = "88.948.82"
$ws.Cells.Item(17, 5).Value = "  +3.30%  "
$ws.Cells.Item(18, 4).Value = "3.215.71"
$ws.Cells.Item(18, 5).Value = "  -0.93%  "
$ws.Cells.Item(19, 4).Value = "'3.22"
$ws.Cells.Item(19, 5).Value = "  +8.03%  "
$ws.Cells.Item(20, 4).Value = "'13.89"
$ws.Cells.Item(20, 5).Value = "  -1.01%  "
$ws.Cells.Item(21, 4).Value = "'421.24"
$ws.Cells.Item(21, 5).Value = "  -1.37%  "
$ws.Cells.Item(22, 4).Value = "'8.74"
$ws.Cells.Item(22, 5).Value = "  -2.72%  "
$ws.Cells.Item(23, 4).Value = "'5.18"
$ws.Cells.Item(23, 5).Value = "  -1.98%  "
$ws.Cells.Item(24, 4).Value = "'0.0000183"
$ws.Cells.Item(24, 5).Value = "  +41.70%  "
$ws.Cells.Item(25, 4).Value = "'5.38"
$ws.Cells.Item(25, 5).Value = "  +4.68%  "
$ws.Cells.Item(26, 4).Value = "'12.26"
$ws.Cells.Item(26, 5).Value = "  -1.50%  "
$ws.Cells.Item(27, 5).Value = "  -2.56%  "
$ws.Cells.Item(28, 4).Value = "'75.06"
$ws.Cells.Item(28, 5).Value = "  -0.94%  "
$ws.Cells.Item(29, 5).Value = "  -0.16%  "
$ws.Cells.Item(30, 4).Value = "'0.171"
$ws.Cells.Item(30, 5).Value = "  -1.26%  "
$ws.Cells.Item(31, 5).Value = "  +0.01%  "
$ws.Cells.Item(32, 4).Value = "'562.32"
$ws.Cells.Item(32, 5).Value = "  +2.74%  "
$ws.Cells.Item(33, 4).Value = "'8.43"
$ws.Cells.Item(33, 5).Value = "  -4.28%  "
$ws.Cells.Item(34, 4).Value = "'7.13"
$ws.Cells.Item(34, 5).Value = "  +5.41%  "
$ws.Cells.Item(35, 5).Value = "  -6.34%  "
$ws.Cells.Item(36, 4).Value = "'1.89"
$ws.Cells.Item(36, 5).Value = "  -2.75%  "
$ws.Cells.Item(37, 2).Value = "dogwifhat"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(37, 4).Value = "'3.34"
$ws.Cells.Item(37, 5).Value = "  +14.00%  "
$ws.Cells.Item(38, 4).Value = "'22.20"
$ws.Cells.Item(38, 5).Value = "  -0.78%  "
$ws.Cells.Item(39, 2).Value = "Kaspa"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(39, 4).Value = "'0.132"
$ws.Cells.Item(39, 5).Value = "  -2.94%  "
$ws.Cells.Item(40, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(40, 4).Value = "'21.87"
$ws.Cells.Item(40, 5).Value = "  +1.04%  "
$ws.Cells.Item(41, 4).Value = "'0.994"
$ws.Cells.Item(41, 5).Value = "  -0.59%  "
$ws.Cells.Item(42, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(42, 4).Value = "'0.386"
$ws.Cells.Item(42, 5).Value = "  -1.28%  "
$ws.Cells.Item(43, 2).Value = "USDe"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(43, 4).Value = "'1.00"
$ws.Cells.Item(43, 5).Value = "  +0.09%  "
$ws.Cells.Item(44, 2).Value = "Stacks"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(44, 4).Value = "'1.95"
$ws.Cells.Item(44, 5).Value = "  -2.48%  "
$ws.Cells.Item(45, 4).Value = "'150.32"
$ws.Cells.Item(45, 5).Value = "  -4.89%  "
$ws.Cells.Item(46, 4).Value = "'180.46"
$ws.Cells.Item(46, 5).Value = "  +1.61%  "
$ws.Cells.Item(47, 4).Value = "'43.56"
$ws.Cells.Item(47, 5).Value = "  -1.61%  "
$ws.Cells.Item(48, 4).Value = "'0.127"
$ws.Cells.Item(48, 5).Value = "  +8.08%  "
$ws.Cells.Item(49, 5).Value = "  -3.24%  "
$ws.Cells.Item(50, 2).Value = "ARBITRUM"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(50, 4).Value = "'0.612"
$ws.Cells.Item(50, 5).Value = "  -1.06%  "
$ws.Cells.Item(51, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(51, 4).Value = "'24.70"
$ws.Cells.Item(51, 5).Value = "  +1.96%  "
